$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.468.31'
$ws.Range('E2').Value = '  -3.85%  '
$ws.Range('D3').Value = '2.505.02'
$ws.Range('E3').Value = '  -4.93%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '576.59'
$ws.Range('D6').Value = '166.14'
$ws.Range('E6').Value = '  -4.51%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '0.517'
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '2.506.50'
$ws.Range('D10').Value = '0.159'
$ws.Range('E10').Value = '  -6.39%  '
$ws.Range('E11').Value = '  -1.48%  '
$ws.Range('E12').Value = '  -2.70%  '
$ws.Range('D13').Value = '4.86'
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').Value = '2.965.58'
$ws.Range('E14').Value = '  -4.81%  '
$ws.Range('D15').Value = '69.357.65'
$ws.Range('E15').Value = '  -3.89%  '
$ws.Range('E16').Value = '  -6.12%  '
$ws.Range('D17').Value = '24.83'
$ws.Range('E17').Value = '  -3.61%  '
$ws.Range('D18').Value = '2.516.96'
$ws.Range('E18').Value = '  -3.78%  '
$ws.Range('D19').Value = '7.81'
$ws.Range('E19').Value = '  -7.71%  '
$ws.Range('E20').Value = '  -6.05%  '
$ws.Range('D21').Value = '347.83'
$ws.Range('E21').Value = '  -6.95%  '
$ws.Range('D22').Value = '3.95'
$ws.Range('E22').Value = '  -3.63%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('E24').Value = '  -5.61%  '
$ws.Range('D25').Value = '68.55'
$ws.Range('E25').Value = '  -3.01%  '
$ws.Range('E26').Value = '  -5.93%  '
$ws.Range('E27').Value = '  -6.85%  '
$ws.Range('E28').Value = '  -5.00%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0900'
$ws.Range('E29').Value = '  -4.95%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '7.87'
$ws.Range('E30').Value = '  -0.58%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '461.88'
$ws.Range('E31').Value = '  -6.75%  '
$ws.Range('E32').Value = '  -1.00%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.75'
$ws.Range('E33').Value = '  -2.12%  '
$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.117'
$ws.Range('E35').Value = '  +2.40%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '154.19'
$ws.Range('E36').Value = '  -5.44%  '
$ws.Range('B37').Value = 'WhiteBITCoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D37').Value = '18.95'
$ws.Range('E37').Value = '  +0.61%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = '18.39'
$ws.Range('E38').Value = '  -4.12%  '
$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').Value = '4.73'
$ws.Range('E40').Value = '  -2.80%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').Value = '0.316'
$ws.Range('E41').Value = '  -2.63%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '1.60'
$ws.Range('E42').Value = '  -6.91%  '
$ws.Range('B43').Value = 'ImmutableX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D43').Value = '1.15'
$ws.Range('E43').Value = '  -14.51%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '38.11'
$ws.Range('E44').Value = '  -2.38%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '2.29'
$ws.Range('E45').Value = '  -10.31%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '143.30'
$ws.Range('E46').Value = '  -5.65%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '0.527'
$ws.Range('E47').Value = '  -2.97%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').Value = '3.50'
$ws.Range('E48').Value = '  -3.88%  '
$ws.Range('B49').Value = 'Optimism'
$ws.Range('C49').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D49').Value = '1.59'
$ws.Range('E49').Value = '  -4.83%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0730'
$ws.Range('E50').Value = '  -1.55%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.583'
$ws.Range('E51').Value = '  -2.77%  '
